# The "Translation" sheet (sheet2) lists UI texts. This change removes the
# first data row (the old "SingleUseId15" entry for the wildcard "<" text)
# which shifts all subsequent rows up by one, and appends two brand new
# rows ("SingleUseId25" and "SingleUseId26") at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Remove the former row 11 (SingleUseId15); everything below shifts up.
$ws.Rows.Item(11).Delete() | Out-Null

# Append new row: SingleUseId25 / Default / Center / LTR / ">"
$ws.Cells.Item(20, 2).Value = "SingleUseId25"
$ws.Cells.Item(20, 3).Value = "Default"
$ws.Cells.Item(20, 4).Value = "Center"
$ws.Cells.Item(20, 5).Value = "LTR"
$ws.Cells.Item(20, 6).Value = ">"

# Append new row: SingleUseId26 / Default / Left / LTR / "New Text"
$ws.Cells.Item(21, 2).Value = "SingleUseId26"
$ws.Cells.Item(21, 3).Value = "Default"
$ws.Cells.Item(21, 4).Value = "Left"
$ws.Cells.Item(21, 5).Value = "LTR"
$ws.Cells.Item(21, 6).Value = "New Text"
